$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.470.08"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.647.89"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9998"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "300.25"
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3803"
$ws.Range("E7").Value = "  -0.91%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "50.54"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3509"
$ws.Range("E9").Value = "  -2.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.223"
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08079"
$ws.Range("E11").Value = "  -0.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9997"
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.09"
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.319"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.267"
$ws.Range("E15").Value = "  -2.23%  "
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").Value = "1.647.25"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.00"
$ws.Range("E18").Value = "  -2.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06969"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.629"
$ws.Range("E20").Value = "  -1.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.46"
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.48"
$ws.Range("E23").Value = "  -0.94%  "
$ws.Range("D24").Value = "23.472.36"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.426"
$ws.Range("E25").Value = "  -2.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.991"
$ws.Range("E26").Value = "  -1.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.05"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.176"
$ws.Range("E29").Value = "  -1.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.02"
$ws.Range("E30").Value = "  -1.38%  "
$ws.Range("D31").Value = "1.841.70"
$ws.Range("E31").Value = "  +0.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.865"
$ws.Range("E32").Value = "  -2.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.133"
$ws.Range("E33").Value = "  -5.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.20"
$ws.Range("E34").Value = "  -8.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9905"
$ws.Range("E35").Value = "  -5.92%  "
$ws.Range("E36").Value = "  -2.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08754"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.921"
$ws.Range("E38").Value = "  -2.15%  "
$ws.Range("E39").Value = "  -2.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06795"
$ws.Range("E40").Value = "  -2.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.83"
$ws.Range("E41").Value = "  -1.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6865"
$ws.Range("E42").Value = "  -1.41%  "
$ws.Range("E43").Value = "  -3.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.53"
$ws.Range("E44").Value = "  -2.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9988"
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6384"
$ws.Range("E46").Value = "  -1.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.251"
$ws.Range("E47").Value = "  -1.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.915"
$ws.Range("E48").Value = "  -0.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07683"
$ws.Range("E49").Value = "  -2.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "127.03"
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("E51").Value = "  +2.48%  "
